# Apply updated cryptocurrency price values to column D (Price)
# Cells are stored as text in the original workbook, so we force a text
# number format before assigning the new value to avoid Excel auto-converting
# the numeric-looking string into a true number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.57"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "24.17"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.304"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05736"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.480"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.129"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8171"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8670"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1377"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07011"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02916"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09393"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.743"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001539"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04712"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005970"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006160"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001243"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.003860"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00008800"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.149"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3174"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1320"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1359"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03715"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006390"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007856"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005276"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.3500"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002140"
